$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed in the source data repull
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -4
$ws.Range("F17").Value = -6
$ws.Range("F18").Value = -4
$ws.Range("F20").Value = -4
$ws.Range("F27").Value = -2
$ws.Range("F29").Value = -7
$ws.Range("F33").Value = -1
